$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 43566
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 3.5

$ws.Range("A7").Value = 43567
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = 8

$ws.Range("K3").Select()
